# Update the "manualStatus" (column I) values from plain numbers to
# bracketed manual-status codes (text), per commit "update to manual
# status column;".
#
# Rows 2-13,15,16,18 : 128   -> "[128]"
# Rows 14,17          : 1128 -> "[1,128]"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 9).Value = "[128]"
$ws.Cells.Item(3, 9).Value = "[128]"
$ws.Cells.Item(4, 9).Value = "[128]"
$ws.Cells.Item(5, 9).Value = "[128]"
$ws.Cells.Item(6, 9).Value = "[128]"
$ws.Cells.Item(7, 9).Value = "[128]"
$ws.Cells.Item(8, 9).Value = "[128]"
$ws.Cells.Item(9, 9).Value = "[128]"
$ws.Cells.Item(10, 9).Value = "[128]"
$ws.Cells.Item(11, 9).Value = "[128]"
$ws.Cells.Item(12, 9).Value = "[128]"
$ws.Cells.Item(13, 9).Value = "[128]"
$ws.Cells.Item(14, 9).Value = "[1,128]"
$ws.Cells.Item(15, 9).Value = "[128]"
$ws.Cells.Item(16, 9).Value = "[128]"
$ws.Cells.Item(17, 9).Value = "[1,128]"
$ws.Cells.Item(18, 9).Value = "[128]"

# The row heights auto-shrank slightly (15 -> 13.8) for all data rows
# except row 14, which kept its original height.
$ws.Rows.Item(3).RowHeight = 13.8
$ws.Rows.Item(4).RowHeight = 13.8
$ws.Rows.Item(5).RowHeight = 13.8
$ws.Rows.Item(6).RowHeight = 13.8
$ws.Rows.Item(7).RowHeight = 13.8
$ws.Rows.Item(8).RowHeight = 13.8
$ws.Rows.Item(9).RowHeight = 13.8
$ws.Rows.Item(10).RowHeight = 13.8
$ws.Rows.Item(11).RowHeight = 13.8
$ws.Rows.Item(12).RowHeight = 13.8
$ws.Rows.Item(13).RowHeight = 13.8
$ws.Rows.Item(15).RowHeight = 13.8
$ws.Rows.Item(16).RowHeight = 13.8
$ws.Rows.Item(17).RowHeight = 13.8
$ws.Rows.Item(18).RowHeight = 13.8

# Selection moved to I18.
$ws.Range("I18").Select() | Out-Null
